$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1393
$ws.Range("I19").Value = 812.75
$ws.Range("J19").Value = 1724.5714
$ws.Range("K19").Value = 812.75
$ws.Range("L19").Value = 1724.5714
$ws.Range("M19").Value = -637.75
$ws.Range("N19").Value = -2074.5714
$ws.Range("H32").Value = 4633
$ws.Range("I32").Value = 1999.5
$ws.Range("J32").Value = 5949.75
$ws.Range("K32").Value = 1999.5
$ws.Range("L32").Value = 5949.75
$ws.Range("M32").Value = -1673.5
$ws.Range("N32").Value = -6601.75
$ws.Range("H33").Value = 810.6316
$ws.Range("I33").Value = 294.53845
$ws.Range("K33").Value = 294.53845
$ws.Range("M33").Value = -65.53845000000001
$ws.Range("H40").Value = 3196
$ws.Range("J40").Value = 3236.818
$ws.Range("L40").Value = 3236.818
$ws.Range("N40").Value = -3586.818
$ws.Range("H41").Value = 398.41666
$ws.Range("I41").Value = 679.8
$ws.Range("J41").Value = 197.42857
$ws.Range("K41").Value = 679.8
$ws.Range("L41").Value = 197.42857
$ws.Range("M41").Value = -239.8
$ws.Range("N41").Value = -1077.42857
$ws.Range("H69").Value = 142863220
$ws.Range("I69").Value = 5837.3335
$ws.Range("J69").Value = 250006260
$ws.Range("K69").Value = 17512.0005
$ws.Range("L69").Value = 750018780
$ws.Range("M69").Value = -16638.0005
$ws.Range("N69").Value = -750020528
$ws.Range("H72").Value = 142863220
$ws.Range("I72").Value = 5837.3335
$ws.Range("J72").Value = 250006260
$ws.Range("K72").Value = 52536.0015
$ws.Range("L72").Value = 2250056340
$ws.Range("M72").Value = -48168.0015
$ws.Range("N72").Value = -2250065076
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").Value = ""
$ws.Range("H76").Value = 3554.2856
$ws.Range("I76").Value = 3608.75
$ws.Range("J76").Value = 3481.6667
$ws.Range("K76").Value = 3608.75
$ws.Range("L76").Value = 3481.6667
$ws.Range("M76").Value = -3293.75
$ws.Range("N76").Value = -4111.6667
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").Value = ""
$ws.Range("H79").Value = 3554.2856
$ws.Range("I79").Value = 3608.75
$ws.Range("J79").Value = 3481.6667
$ws.Range("K79").Value = 3608.75
$ws.Range("L79").Value = 3481.6667
$ws.Range("M79").Value = -2516.75
$ws.Range("N79").Value = -5665.6667
$ws.Range("H80").Value = 2091.125
$ws.Range("J80").Value = 2987.5
$ws.Range("L80").Value = 8962.5
$ws.Range("N80").Value = -10958.5
$ws.Range("H83").Value = 2091.125
$ws.Range("J83").Value = 2987.5
$ws.Range("L83").Value = 26887.5
$ws.Range("N83").Value = -36871.5
$ws.Range("H86").Value = 4095.1
$ws.Range("J86").Value = 4650.3335
$ws.Range("L86").Value = 4650.3335
$ws.Range("N86").Value = -6896.3335
$ws.Range("H88").Value = 5869
$ws.Range("J88").Value = 5869
$ws.Range("L88").Value = 5869
$ws.Range("N88").Value = -6681
$ws.Range("H89").Value = 4095.1
$ws.Range("J89").Value = 4650.3335
$ws.Range("L89").Value = 23251.6675
$ws.Range("N89").Value = -34483.6675
$ws.Range("H91").Value = 5869
$ws.Range("J91").Value = 5869
$ws.Range("L91").Value = 5869
$ws.Range("N91").Value = -8677
$ws.Range("H92").Value = 1071.7826
$ws.Range("I92").Value = 856.4286
$ws.Range("K92").Value = 856.4286
$ws.Range("M92").Value = 391.5714
$ws.Range("H96").Value = 834555.75
$ws.Range("I96").Value = 1429167
$ws.Range("K96").Value = 4287501
$ws.Range("M96").Value = -4286128
$ws.Range("H106").Value = 8573.519
$ws.Range("I106").Value = 7867.857
$ws.Range("K106").Value = 7867.857
$ws.Range("M106").Value = -7236.857
$ws.Range("H111").Value = 1301.5834
$ws.Range("I111").Value = 2428.5
$ws.Range("J111").Value = 738.125
$ws.Range("K111").Value = 7285.5
$ws.Range("L111").Value = 2214.375
$ws.Range("M111").Value = -4218.5
$ws.Range("N111").Value = -8348.375
$ws.Range("H125").Value = 35001.168
$ws.Range("J125").Value = 41602.6
$ws.Range("L125").Value = 374423.4
$ws.Range("N125").Value = -379343.4
$ws.Range("H137").Value = 3145.9333
$ws.Range("I137").Value = 3455.8
$ws.Range("K137").Value = 10367.4
$ws.Range("M137").Value = -7817.400000000001
$ws.Range("H141").Value = 3223.6897
$ws.Range("I141").Value = 3055.074
$ws.Range("J141").Value = 5500
$ws.Range("K141").Value = 9165.222
$ws.Range("L141").Value = 16500
$ws.Range("M141").Value = -3985.222
$ws.Range("N141").Value = -26860
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2994.6833
$ws.Range("I32").Value = 2553.4822
$ws.Range("J32").Value = 9171.5
$ws.Range("K32").Value = 2553.4822
$ws.Range("L32").Value = 9171.5
$ws.Range("M32").Value = -2266.4822
$ws.Range("N32").Value = -9745.5
$ws.Range("H36").Value = 2854
$ws.Range("I36").Value = 2854
$ws.Range("K36").Value = 2854
$ws.Range("M36").Value = -2508
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("N40").Value = ""
$ws.Range("H63").Value = 9999
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").Value = ""
$ws.Range("H66").Value = 9999
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").Value = ""
$ws.Range("H97").Value = 721.5714
$ws.Range("I97").Value = 559.44446
$ws.Range("K97").Value = 559.44446
$ws.Range("M97").Value = -63.44446000000005
$ws.Range("H110").Value = 2526.125
$ws.Range("I110").Value = 1639.2
$ws.Range("K110").Value = 1639.2
$ws.Range("M110").Value = 405.8
$ws.Range("H122").Value = 4413.5
$ws.Range("I122").Value = 4307.75
$ws.Range("K122").Value = 12923.25
$ws.Range("M122").Value = -10473.25
$ws.Range("H132").Value = 2547.8262
$ws.Range("I132").Value = 2646.762
$ws.Range("J132").Value = 1509
$ws.Range("K132").Value = 7940.286
$ws.Range("L132").Value = 4527
$ws.Range("M132").Value = -5410.286
$ws.Range("N132").Value = -9587
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").Value = ""
$ws.Range("H20").Value = 1278.0834
$ws.Range("I20").Value = 872.4211
$ws.Range("K20").Value = 872.4211
$ws.Range("M20").Value = -625.4211
$ws.Range("H75").Value = 15566.167
$ws.Range("I75").Value = 15566.167
$ws.Range("K75").Value = 15566.167
$ws.Range("M75").Value = -14630.167
$ws.Range("H78").Value = 15566.167
$ws.Range("I78").Value = 15566.167
$ws.Range("K78").Value = 46698.501
$ws.Range("M78").Value = -42018.501
$ws.Range("H86").Value = 25952248
$ws.Range("I86").Value = 1996.5
$ws.Range("K86").Value = 1996.5
$ws.Range("M86").Value = -873.5
$ws.Range("H89").Value = 25952248
$ws.Range("I89").Value = 1996.5
$ws.Range("K89").Value = 9982.5
$ws.Range("M89").Value = -4366.5
$ws.Range("H99").Value = 3693.5715
$ws.Range("I99").Value = 1526.25
$ws.Range("K99").Value = 1526.25
$ws.Range("M99").Value = -28.25
$ws.Range("H105").Value = 3700.5454
$ws.Range("I105").Value = 3400.0476
$ws.Range("K105").Value = 3400.0476
$ws.Range("M105").Value = -1653.0476
$ws.Range("H134").Value = 2586.2856
$ws.Range("I134").Value = 2584.0952
$ws.Range("K134").Value = 7752.285600000001
$ws.Range("M134").Value = -5217.285600000001
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1647.3846
$ws.Range("I16").Value = 1224.4286
$ws.Range("K16").Value = 1224.4286
$ws.Range("M16").Value = -937.4286
$ws.Range("H31").Value = 4571.408
$ws.Range("I31").Value = 2420
$ws.Range("K31").Value = 2420
$ws.Range("M31").Value = -2125
$ws.Range("H34").Value = 4571.408
$ws.Range("I34").Value = 2420
$ws.Range("K34").Value = 2420
$ws.Range("M34").Value = -2218
$ws.Range("H62").Value = 7315.909
$ws.Range("I62").Value = 5408.6
$ws.Range("K62").Value = 5408.6
$ws.Range("M62").Value = -4784.6
$ws.Range("H65").Value = 7315.909
$ws.Range("I65").Value = 5408.6
$ws.Range("K65").Value = 27043
$ws.Range("M65").Value = -23923
$ws.Range("H86").Value = 5855.8887
$ws.Range("J86").Value = 5882.3335
$ws.Range("L86").Value = 5882.3335
$ws.Range("N86").Value = -8128.3335
$ws.Range("H89").Value = 5855.8887
$ws.Range("J89").Value = 5882.3335
$ws.Range("L89").Value = 29411.6675
$ws.Range("N89").Value = -40643.6675
$ws.Range("H99").Value = 3120
$ws.Range("I99").Value = 3120
$ws.Range("K99").Value = 3120
$ws.Range("M99").Value = -1622
$ws.Range("H105").Value = 2977.3
$ws.Range("I105").Value = 3710.2856
$ws.Range("K105").Value = 3710.2856
$ws.Range("M105").Value = -1963.2856
$ws.Range("H113").Value = 1647.3846
$ws.Range("I113").Value = 1224.4286
$ws.Range("K113").Value = 1224.4286
$ws.Range("M113").Value = 945.5714
$ws.Range("H126").Value = 3120
$ws.Range("I126").Value = 3120
$ws.Range("K126").Value = 9360
$ws.Range("M126").Value = -6890
$ws.Range("H132").Value = 2222.8125
$ws.Range("I132").Value = 2171
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 6513
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -3983
$ws.Range("N132").Value = -14060
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 13742477
$ws.Range("I4").Value = 365162.44
$ws.Range("J4").Value = 101303090
$ws.Range("K4").Value = 1095487.32
$ws.Range("L4").Value = 303909270
$ws.Range("M4").Value = -1095375.32
$ws.Range("N4").Value = -303909494
$ws.Range("H107").Value = 1181.3914
$ws.Range("I107").Value = 3552.8333
$ws.Range("K107").Value = 10658.4999
$ws.Range("M107").Value = -8738.499899999999
$ws.Range("H113").Value = 2723.2666
$ws.Range("J113").Value = 3183.25
$ws.Range("L113").Value = 9549.75
$ws.Range("N113").Value = -13889.75
$ws.Range("H114").Value = 1749.8
$ws.Range("I114").Value = 937.25
$ws.Range("J114").Value = 5000
$ws.Range("K114").Value = 2811.75
$ws.Range("L114").Value = 15000
$ws.Range("M114").Value = 442.25
$ws.Range("N114").Value = -21508
$ws.Range("H122").Value = 505.625
$ws.Range("J122").Value = 696
$ws.Range("L122").Value = 6264
$ws.Range("N122").Value = -11164
$ws.Range("H130").Value = 4993.9
$ws.Range("I130").Value = 5051.857
$ws.Range("K130").Value = 15155.571
$ws.Range("M130").Value = -10135.571
$ws.Range("H138").Value = 3991
$ws.Range("J138").Value = 3991
$ws.Range("L138").Value = 11973
$ws.Range("N138").Value = -22253
$ws.Range("H140").Value = 1293.625
$ws.Range("I140").Value = 1049.8572
$ws.Range("K140").Value = 3149.5716
$ws.Range("M140").Value = 2030.4284
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 36.175
$ws.Range("J2").Value = 39
$ws.Range("L2").Value = 39
$ws.Range("N2").Value = -265
$ws.Range("H3").Value = 3306.5715
$ws.Range("I3").Value = 2524.3333
$ws.Range("J3").Value = 8000
$ws.Range("K3").Value = 2524.3333
$ws.Range("L3").Value = 8000
$ws.Range("M3").Value = -2408.3333
$ws.Range("N3").Value = -8232
$ws.Range("H80").Value = 3139.2
$ws.Range("I80").Value = 2296.625
$ws.Range("J80").Value = 3700.9167
$ws.Range("K80").Value = 2296.625
$ws.Range("L80").Value = 3700.9167
$ws.Range("M80").Value = -1298.625
$ws.Range("N80").Value = -5696.9167
$ws.Range("H83").Value = 3139.2
$ws.Range("I83").Value = 2296.625
$ws.Range("J83").Value = 3700.9167
$ws.Range("K83").Value = 11483.125
$ws.Range("L83").Value = 18504.5835
$ws.Range("M83").Value = -6491.125
$ws.Range("N83").Value = -28488.5835
$ws.Range("H97").Value = 876.1667
$ws.Range("I97").Value = 651.73334
$ws.Range("K97").Value = 651.73334
$ws.Range("M97").Value = -155.73334
$ws.Range("H102").Value = 29999.6
$ws.Range("I102").Value = 29999
$ws.Range("K102").Value = 29999
$ws.Range("M102").Value = -28377
$ws.Range("H122").Value = 9493.799999999999
$ws.Range("I122").Value = 11360.3
$ws.Range("J122").Value = 5760.8
$ws.Range("K122").Value = 34080.89999999999
$ws.Range("L122").Value = 17282.4
$ws.Range("M122").Value = -31630.89999999999
$ws.Range("N122").Value = -22182.4
$ws.Range("H132").Value = 2240.9688
$ws.Range("I132").Value = 1681.6957
$ws.Range("J132").Value = 3670.2222
$ws.Range("K132").Value = 5045.0871
$ws.Range("L132").Value = 11010.6666
$ws.Range("M132").Value = -2515.0871
$ws.Range("N132").Value = -16070.6666
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2522.1667
$ws.Range("I7").Value = 2471
$ws.Range("K7").Value = 2471
$ws.Range("M7").Value = -2359
$ws.Range("H40").Value = 4893.9287
$ws.Range("I40").Value = 3807.0667
$ws.Range("K40").Value = 3807.0667
$ws.Range("M40").Value = -3671.0667
$ws.Range("H46").Value = 3165.6875
$ws.Range("J46").Value = 3425.0715
$ws.Range("L46").Value = 3425.0715
$ws.Range("N46").Value = -3801.0715
$ws.Range("H61").Value = 3978.3845
$ws.Range("I61").Value = 3175.32
$ws.Range("K61").Value = 3175.32
$ws.Range("M61").Value = -2973.32
$ws.Range("H68").Value = 6853.3335
$ws.Range("I68").Value = 4539.4
$ws.Range("J68").Value = 8010.3
$ws.Range("K68").Value = 4539.4
$ws.Range("L68").Value = 8010.3
$ws.Range("M68").Value = -3790.4
$ws.Range("N68").Value = -9508.299999999999
$ws.Range("H71").Value = 6853.3335
$ws.Range("I71").Value = 4539.4
$ws.Range("J71").Value = 8010.3
$ws.Range("K71").Value = 22697
$ws.Range("L71").Value = 40051.5
$ws.Range("M71").Value = -18953
$ws.Range("N71").Value = -47539.5
$ws.Range("H93").Value = 5716.6665
$ws.Range("I93").Value = 1100
$ws.Range("J93").Value = 6640
$ws.Range("K93").Value = 1100
$ws.Range("L93").Value = 6640
$ws.Range("M93").Value = 148
$ws.Range("N93").Value = -9136
$ws.Range("H94").Value = 60917.668
$ws.Range("J94").Value = 51226.5
$ws.Range("L94").Value = 51226.5
$ws.Range("N94").Value = -52578.5
$ws.Range("H113").Value = 3978.3845
$ws.Range("I113").Value = 3175.32
$ws.Range("K113").Value = 3175.32
$ws.Range("M113").Value = -1005.32
$ws.Range("H122").Value = 11208.125
$ws.Range("I122").Value = 7443.2
$ws.Range("K122").Value = 22329.6
$ws.Range("M122").Value = -19879.6
$ws.Range("H126").Value = 2522.1667
$ws.Range("I126").Value = 2471
$ws.Range("K126").Value = 7413
$ws.Range("M126").Value = -4943
$ws.Range("H132").Value = 2972.9412
$ws.Range("I132").Value = 2369.3076
$ws.Range("J132").Value = 4934.75
$ws.Range("K132").Value = 7107.9228
$ws.Range("L132").Value = 14804.25
$ws.Range("M132").Value = -4577.9228
$ws.Range("N132").Value = -19864.25
$ws.Range("H136").Value = 17109.5
$ws.Range("I136").Value = 1278
$ws.Range("K136").Value = 3834
$ws.Range("M136").Value = -1284
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 3926.389
$ws.Range("J107").Value = 4860.75
$ws.Range("L107").Value = 14582.25
$ws.Range("N107").Value = -18422.25
$ws.Range("H126").Value = 2790.4375
$ws.Range("I126").Value = 2513.84
$ws.Range("K126").Value = 7541.52
$ws.Range("M126").Value = -5071.52
$ws.Range("H132").Value = 2845.6177
$ws.Range("I132").Value = 1265.4166
$ws.Range("J132").Value = 3707.5454
$ws.Range("K132").Value = 3796.2498
$ws.Range("L132").Value = 11122.6362
$ws.Range("M132").Value = -1266.2498
$ws.Range("N132").Value = -16182.6362
$ws.Range("H136").Value = 1983.0454
$ws.Range("I136").Value = 1574.9333
$ws.Range("K136").Value = 4724.7999
$ws.Range("M136").Value = -2174.7999
